# "10Th - MB for single stock and added new group"
#
# This report table has one row per analyst/broker (column A) and one
# column per trading day (most-recent day first, starting at column B).
# Two new days (Jun_26, Jun_27) are being published, so three new date
# columns are inserted in front of the existing data (the existing Jun_17
# column and everything to its right shifts three columns to the right),
# and two new broker rows (Benchmark, Evercore ISI) are appended at the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new date columns before the existing first data column (B) ---
# This shifts old B->E, C->F, D->G, E->H, carrying their values/styles
# along for the ride (including the highlighted note cell, old C5 -> F5).
$ws.Columns("B:D").Insert()

# Re-apply the (lost-on-insert) custom width to every data column so they
# all keep the original 8-character width.
$ws.Columns("C:H").ColumnWidth = 7.1667

# --- New header row values for the freshly inserted columns ---
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Fill in the new columns for every existing analyst row with "UN" ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- Two new broker/analyst rows appended at the bottom ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
